$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values
$ws.Range("B2").Value = 522.40334535812497
$ws.Range("C2").Value = 437.90090132812503
$ws.Range("D2").Value = 522.40334535812497
$ws.Range("E2").Value = 437.90090132812503

# Row 3 data values
$ws.Range("B3").Value = 526.34982737125006
$ws.Range("C3").Value = 432.99391268812496
$ws.Range("D3").Value = 526.34982737125006
$ws.Range("E3").Value = 432.99391268812496

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
